$wb = $excel.ActiveWorkbook

# --- Flight Mission Cycle sheet: now totals duration by summing each setting ---
$flight = $wb.Worksheets.Item("Flight Mission Cycle")

# Row 2: Typing duration, summed from the Typing sheet's per-axis durations (B3:D3)
$flight.Cells.Item(2, 1).Value = "Typing"
$flight.Cells.Item(2, 2).Formula = "=SUM(Typing!B3:D3)"

# Row 3: Light switch duration
$flight.Cells.Item(3, 1).Value = "Light switch"
$flight.Cells.Item(3, 2).Value = 20

# Row 4: Piano duration
$flight.Cells.Item(4, 1).Value = "Piano"
$flight.Cells.Item(4, 2).Value = 30

# Row 5: Writing duration (previously in row 2)
$flight.Cells.Item(5, 1).Value = "Writing"
$flight.Cells.Item(5, 2).Value = 50

$flight.Range("K18").Select()

# --- Typing sheet: set Period value ---
$typing = $wb.Worksheets.Item("Typing")
$typing.Cells.Item(6, 2).Value = 10
$typing.Range("I19").Select()

# --- Light switch sheet: set Max_RoM and Duration values ---
$lightSwitch = $wb.Worksheets.Item("Light switch")
$lightSwitch.Cells.Item(2, 2).Value = 10
$lightSwitch.Cells.Item(3, 2).Value = 10
$lightSwitch.Range("F20").Select()

# --- Writing sheet: just selection moved ---
$writing = $wb.Worksheets.Item("Writing")
$writing.Range("H26").Select()

# --- Piano sheet: set Max_RoM value; ends up the active sheet/tab ---
$piano = $wb.Worksheets.Item("Piano")
$piano.Cells.Item(4, 2).Value = 30
$piano.Activate()
$piano.Range("U5").Select()
